# Update "想去人数" (column F) counts for the latest data refresh
# (commit: "Update gh-pages to output generated at 456a3b4")

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1893
$ws.Cells.Item(8, 6).Value = 943
$ws.Cells.Item(10, 6).Value = 1291
$ws.Cells.Item(11, 6).Value = 1572
$ws.Cells.Item(13, 6).Value = 1571
$ws.Cells.Item(14, 6).Value = 350
$ws.Cells.Item(17, 6).Value = 1138
$ws.Cells.Item(21, 6).Value = 1850
$ws.Cells.Item(24, 6).Value = 1010
$ws.Cells.Item(26, 6).Value = 1275
$ws.Cells.Item(30, 6).Value = 1199
$ws.Cells.Item(35, 6).Value = 290
$ws.Cells.Item(37, 6).Value = 896
$ws.Cells.Item(39, 6).Value = 1706
$ws.Cells.Item(40, 6).Value = 15
$ws.Cells.Item(42, 6).Value = 4
$ws.Cells.Item(43, 6).Value = 2079
$ws.Cells.Item(46, 6).Value = 17

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 271
$ws.Cells.Item(6, 6).Value = 4698
$ws.Cells.Item(10, 6).Value = 785
$ws.Cells.Item(13, 6).Value = 1155
$ws.Cells.Item(15, 6).Value = 768

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 1893
$ws.Cells.Item(3, 6).Value = 271
$ws.Cells.Item(5, 6).Value = 4698
$ws.Cells.Item(6, 6).Value = 785
$ws.Cells.Item(11, 6).Value = 943
$ws.Cells.Item(13, 6).Value = 1291
$ws.Cells.Item(14, 6).Value = 1572
$ws.Cells.Item(16, 6).Value = 1571
$ws.Cells.Item(20, 6).Value = 1138
$ws.Cells.Item(23, 6).Value = 768
$ws.Cells.Item(24, 6).Value = 768
$ws.Cells.Item(25, 6).Value = 1850
$ws.Cells.Item(28, 6).Value = 1010
$ws.Cells.Item(30, 6).Value = 1275
$ws.Cells.Item(34, 6).Value = 1199
$ws.Cells.Item(41, 6).Value = 290
$ws.Cells.Item(42, 6).Value = 896
$ws.Cells.Item(44, 6).Value = 1706
$ws.Cells.Item(45, 6).Value = 15
$ws.Cells.Item(46, 6).Value = 4
$ws.Cells.Item(47, 6).Value = 2079
